# Refitting NCDEs to individual patients (for manuscript figure)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Label" header in column H, styled like the other header cells (row 1)
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Updated (refit) Prediction/Error values for the 100-iteration block (rows 2-8)
$ws.Range("D2").Value = 0.7577676759724691
$ws.Range("E2").Value = 0.7577676759724691

$ws.Range("D3").Value = 0.4858579289840221
$ws.Range("E3").Value = 0.4858579289840221

$ws.Range("D4").Value = 0.5675942854668835
$ws.Range("E4").Value = 0.4324057145331165

$ws.Range("D5").Value = 0.6752220364823475
$ws.Range("E5").Value = 0.3247779635176525

$ws.Range("D6").Value = 0.4777434721531201
$ws.Range("E6").Value = 0.5222565278468799

$ws.Range("D7").Value = 0.4237206033863964
$ws.Range("E7").Value = 0.5762793966136036

$ws.Range("D8").Value = 0.7572762616941298
$ws.Range("E8").Value = 0.2427237383058702
$ws.Range("F8").Value = 0.702509343624115
$ws.Range("G8").Value = 0.5714285714285714

# New "Label" column values (ground-truth diagnosis label per patient) for every data row
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("H15").Value = 1
